$d = $word.ActiveDocument

# Locate the target paragraph: "{m:'Table1'.simpleTable().asStyle('TableauGrille6Couleur-Accentuation3')}"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("{m:'Table1'")) {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# Helper: force a run split boundary at an absolute character position by
# touching a tiny range and toggling (then reverting) a character format -
# this makes the engine re-serialize that run boundary without actually
# changing any visible formatting.
function Split-RunAt($pos) {
    $r = $d.Range($pos, $pos + 1)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# 1) "{m" -> "{" + "m"   (split right after the opening brace)
Split-RunAt ($pStart + 0)

# 2) Insert a space between "Table" and "1" so "Table1" -> "Table 1"
$insPos = $pStart + 9   # position right before the "1" in "Table1"
$insRange = $d.Range($insPos, $insPos)
$insRange.InsertBefore(" ")

# After inserting one character, everything from $insPos onward shifts by +1
# 3) ")}" -> ")" + "}" (split right before the closing brace)
$closeBracePos = $pStart + 72   # position of "}" after the insertion shift
Split-RunAt $closeBracePos

# The very last run (just "}") must end up WITHOUT the <w:rPr><w:lang .../></w:rPr>
# that all the other runs carry. Clear its language formatting so it has no rPr.
$lastCharPos = $target.Range.End - 1
$lastRange = $d.Range($lastCharPos, $target.Range.End)
$lastRange.LanguageID = 0
